$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go headcount) values per the source update.
# The "全部类型" sheet duplicates rows from the other three sheets (merged by date),
# so the same events are updated there too, at different row numbers.

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 57  # was 56
$ws.Range("F3").Value = 1008  # was 1004
$ws.Range("F5").Value = 454  # was 452
$ws.Range("F6").Value = 720  # was 715
$ws.Range("F13").Value = 839  # was 837
$ws.Range("F14").Value = 116  # was 115
$ws.Range("F15").Value = 1986  # was 1984
$ws.Range("F16").Value = 482  # was 479
$ws.Range("F17").Value = 7186  # was 7145
$ws.Range("F18").Value = 527  # was 525
$ws.Range("F20").Value = 55  # was 54

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 55  # was 54

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5481  # was 5474

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 57  # was 56
$ws.Range("F3").Value = 5481  # was 5474
$ws.Range("F7").Value = 1008  # was 1004
$ws.Range("F11").Value = 454  # was 452
$ws.Range("F12").Value = 720  # was 715
$ws.Range("F22").Value = 839  # was 837
$ws.Range("F23").Value = 116  # was 115
$ws.Range("F25").Value = 1986  # was 1984
$ws.Range("F26").Value = 482  # was 479
$ws.Range("F27").Value = 7186  # was 7145
$ws.Range("F28").Value = 55  # was 54
$ws.Range("F29").Value = 527  # was 525
$ws.Range("F31").Value = 55  # was 54
